$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79 - this pushes the existing rows 79..130
# down to 80..131, carrying their content/formatting with them.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly record.
$ws.Range("A79").Value = 3
$ws.Range("B79").Value = "Femacal de La Calera"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 44438
$ws.Range("E79").Value = 5
$ws.Range("F79").Value = 100112001
$ws.Range("G79").Value = "Berenjena"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 70
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 9000
$ws.Range("M79").Value = 9000
$ws.Range("N79").Value = "$/caja 60 unidades"
$ws.Range("O79").Value = "Región de Arica y Parinacota"
$ws.Range("P79").Value = 150
$ws.Range("Q79").Value = 60
$ws.Range("R79").Value = "Hortaliza"
